$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.980.58'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.46%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.418.19'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.67%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '489.54'
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.01'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.64%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.616'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +19.29%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.998'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.12%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.433.29'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.81%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.36'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +10.37%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.100'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.34%  '

$ws.Range("E12").Value = '  -1.51%  '

$ws.Range("E13").Value = '  +1.24%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.829.26'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.97%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '57.083.84'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.49%  '

$ws.Range("E16").Value = '  -3.81%  '

$ws.Range("E17").Value = '  -3.92%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.431.74'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.90%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.70'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.43%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '324.74'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.15%  '

$ws.Range("E21").Value = '  -3.26%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.997'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.93'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.15%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '57.72'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.54%  '

$ws.Range("E25").Value = '  -1.69%  '

$ws.Range("E26").Value = '  -0.22%  '

$ws.Range("E27").Value = '  -2.24%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.520.86'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.72%  '

$ws.Range("E29").Value = '  -5.04%  '

$ws.Range("E30").Value = '  -6.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '151.15'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.62'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.42%  '

$ws.Range("E34").Value = '  -0.62%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.30'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.44%  '

$ws.Range("E36").Value = '  -1.21%  '

$ws.Range("E37").Value = '  -1.47%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.829'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.11%  '

$ws.Range("E39").Value = '  +8.85%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '34.05'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.03%  '

$ws.Range("E41").Value = '  -0.85%  '

$ws.Range("E42").Value = '  -3.53%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '279.54'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.49%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.994'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.598'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.88%  '

$ws.Range("E46").Value = '  -5.85%  '

$ws.Range("E47").Value = '  -0.24%  '

$ws.Range("E48").Value = '  -1.70%  '

$ws.Range("E49").Value = '  -8.78%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.901.41'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.17%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.55'
$ws.Range("D51").Style = "Normal"
